$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2638.3076
$ws.Range("J43").Value = 2654.4546
$ws.Range("L43").Value = 2654.4546
$ws.Range("N43").Value = -2792.4546
$ws.Range("H81").Value = 94664
$ws.Range("J81").Value = 94664
$ws.Range("L81").Value = 94664
$ws.Range("N81").Value = -96660
$ws.Range("H82").Value = 1895.5555
$ws.Range("I82").Value = 1865.7142
$ws.Range("K82").Value = 5597.142599999999
$ws.Range("M82").Value = -5191.142599999999
$ws.Range("H84").Value = 94664
$ws.Range("J84").Value = 94664
$ws.Range("L84").Value = 283992
$ws.Range("N84").Value = -293976
$ws.Range("H85").Value = 1895.5555
$ws.Range("I85").Value = 1865.7142
$ws.Range("K85").Value = 5597.142599999999
$ws.Range("M85").Value = -4193.142599999999
$ws.Range("H112").Value = 5018.6665
$ws.Range("J112").Value = 5108.4253
$ws.Range("L112").Value = 15325.2759
$ws.Range("N112").Value = -17541.2759
$ws.Range("H129").Value = 1512.375
$ws.Range("I129").Value = 442.85715
$ws.Range("J129").Value = 1952.7646
$ws.Range("K129").Value = 1328.57145
$ws.Range("L129").Value = 5858.293799999999
$ws.Range("M129").Value = 3671.42855
$ws.Range("N129").Value = -15858.2938
$ws.Range("H137").Value = 3474544.5
$ws.Range("I137").Value = 4168704
$ws.Range("J137").Value = 3745.75
$ws.Range("K137").Value = 12506112
$ws.Range("L137").Value = 11237.25
$ws.Range("M137").Value = -12503562
$ws.Range("N137").Value = -16337.25
$ws.Range("H138").Value = 5505.839
$ws.Range("I138").Value = 7682.3335
$ws.Range("J138").Value = 4983.48
$ws.Range("K138").Value = 23047.0005
$ws.Range("L138").Value = 14950.44
$ws.Range("M138").Value = -17907.0005
$ws.Range("N138").Value = -25230.44
$ws.Range("H140").Value = 76508.17999999999
$ws.Range("J140").Value = 76508.17999999999
$ws.Range("L140").Value = 76508.17999999999
$ws.Range("N140").Value = -86868.17999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1473.3334
$ws.Range("I2").Value = 1335
$ws.Range("K2").Value = 1335
$ws.Range("M2").Value = -1222
$ws.Range("H32").Value = 22241454
$ws.Range("I32").Value = 27045392
$ws.Range("J32").Value = 23243.75
$ws.Range("K32").Value = 27045392
$ws.Range("L32").Value = 23243.75
$ws.Range("M32").Value = -27045105
$ws.Range("N32").Value = -23817.75
$ws.Range("H61").Value = 7096264
$ws.Range("I61").Value = 16668596
$ws.Range("J61").Value = 5647.5186
$ws.Range("K61").Value = 16668596
$ws.Range("L61").Value = 5647.5186
$ws.Range("M61").Value = -16668384
$ws.Range("N61").Value = -6071.5186
$ws.Range("H116").Value = 1473.3334
$ws.Range("I116").Value = 1335
$ws.Range("K116").Value = 1335
$ws.Range("M116").Value = 959
$ws.Range("H132").Value = 1329294.6
$ws.Range("I132").Value = 2121.5117
$ws.Range("J132").Value = 5133857.5
$ws.Range("K132").Value = 6364.5351
$ws.Range("L132").Value = 15401572.5
$ws.Range("M132").Value = -3834.5351
$ws.Range("N132").Value = -15406632.5
$ws.Range("H136").Value = 7096264
$ws.Range("I136").Value = 16668596
$ws.Range("J136").Value = 5647.5186
$ws.Range("K136").Value = 50005788
$ws.Range("L136").Value = 16942.5558
$ws.Range("M136").Value = -50003238
$ws.Range("N136").Value = -22042.5558

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1473.3334
$ws.Range("I3").Value = 1335
$ws.Range("K3").Value = 1335
$ws.Range("M3").Value = -1221

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7449.9062
$ws.Range("I31").Value = 1067.7778
$ws.Range("J31").Value = 15655.5
$ws.Range("K31").Value = 1067.7778
$ws.Range("L31").Value = 15655.5
$ws.Range("M31").Value = -772.7778000000001
$ws.Range("N31").Value = -16245.5
$ws.Range("H34").Value = 7449.9062
$ws.Range("I34").Value = 1067.7778
$ws.Range("J34").Value = 15655.5
$ws.Range("K34").Value = 1067.7778
$ws.Range("L34").Value = 15655.5
$ws.Range("M34").Value = -865.7778000000001
$ws.Range("N34").Value = -16059.5
$ws.Range("H58").Value = 2622.3333
$ws.Range("I58").Value = 2622.3333
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 2622.3333
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -2419.3333
$ws.Range("N58").ClearContents()
$ws.Range("H132").Value = 35716964
$ws.Range("I132").Value = 83335736
$ws.Range("J132").Value = 2887.25
$ws.Range("K132").Value = 250007208
$ws.Range("L132").Value = 8661.75
$ws.Range("M132").Value = -250004678
$ws.Range("N132").Value = -13721.75
$ws.Range("H134").Value = 5439025
$ws.Range("I134").Value = 6101540.5
$ws.Range("J134").Value = 6396.8
$ws.Range("K134").Value = 18304621.5
$ws.Range("L134").Value = 19190.4
$ws.Range("M134").Value = -18302086.5
$ws.Range("N134").Value = -24260.4
$ws.Range("H136").Value = 2622.3333
$ws.Range("I136").Value = 2622.3333
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 7866.999899999999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -5316.999899999999
$ws.Range("N136").Value = -15316.999899999999
$ws.Range("H141").Value = 106987.46
$ws.Range("J141").Value = 103542.555
$ws.Range("L141").Value = 103542.555
$ws.Range("N141").Value = -113902.555

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 1860
$ws.Range("I48").Value = 1650
$ws.Range("J48").Value = 2000
$ws.Range("K48").Value = 4950
$ws.Range("L48").Value = 6000
$ws.Range("M48").Value = -4700
$ws.Range("N48").Value = -6500
$ws.Range("H56").Value = 4726.6665
$ws.Range("I56").Value = 4726.6665
$ws.Range("K56").Value = 4726.6665
$ws.Range("M56").Value = -4196.6665
$ws.Range("H107").Value = 41668650
$ws.Range("I107").Value = 483.33334
$ws.Range("K107").Value = 1450.00002
$ws.Range("M107").Value = 469.9999800000001
$ws.Range("H113").Value = 1179.5927
$ws.Range("J113").Value = 1550.6923
$ws.Range("L113").Value = 4652.0769
$ws.Range("N113").Value = -8992.0769
$ws.Range("H132").Value = 2682.4285
$ws.Range("I132").Value = 2801
$ws.Range("J132").Value = 2579.6667
$ws.Range("K132").Value = 25209
$ws.Range("L132").Value = 23217.0003
$ws.Range("M132").Value = -22679
$ws.Range("N132").Value = -28277.0003
$ws.Range("H140").Value = 1348.9491
$ws.Range("I140").Value = 903.0571
$ws.Range("K140").Value = 2709.1713
$ws.Range("M140").Value = 2470.8287

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 8213310.5
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
$ws.Range("H132").Value = 38468420
$ws.Range("I132").Value = 76933416
$ws.Range("J132").Value = 3419.6155
$ws.Range("K132").Value = 230800248
$ws.Range("L132").Value = 10258.8465
$ws.Range("M132").Value = -230797718
$ws.Range("N132").Value = -15318.8465

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1620.5
$ws.Range("I16").Value = 1577.2
$ws.Range("J16").Value = 1651.4286
$ws.Range("K16").Value = 1577.2
$ws.Range("L16").Value = 1651.4286
$ws.Range("M16").Value = -1407.2
$ws.Range("N16").Value = -1991.4286
$ws.Range("H40").Value = 4396.0527
$ws.Range("I40").Value = 4043.3333
$ws.Range("J40").Value = 5000.7144
$ws.Range("K40").Value = 4043.3333
$ws.Range("L40").Value = 5000.7144
$ws.Range("M40").Value = -3907.3333
$ws.Range("N40").Value = -5272.7144
$ws.Range("H46").Value = 1434.7273
$ws.Range("I46").Value = 460.33334
$ws.Range("J46").Value = 1800.125
$ws.Range("K46").Value = 460.33334
$ws.Range("L46").Value = 1800.125
$ws.Range("M46").Value = -272.33334
$ws.Range("N46").Value = -2176.125
$ws.Range("H63").Value = 44042.5
$ws.Range("J63").Value = 44042.5
$ws.Range("L63").Value = 44042.5
$ws.Range("N63").Value = -45540.5
$ws.Range("H66").Value = 44042.5
$ws.Range("J66").Value = 44042.5
$ws.Range("L66").Value = 132127.5
$ws.Range("N66").Value = -139615.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H140").Value = 50702.43
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 50702.43
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 50702.43
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -61062.43
